$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing rows (including the old header row 1)
# shift down by one.
$ws.Rows.Item(1).Insert()

# New row 1: a note spanning A1:D1, merged, italic red text.
$ws.Range("A1").Value = "Note: The date header (Row 2) supports: '2023 Annual', '2023 Q1', '2023-01'"
$ws.Range("A1:D1").Merge()
$ws.Range("A1").Font.Italic = $true
$ws.Range("A1").Font.Color = 255

# Row 2 (former row 1) is the account/date header row - update the year
# labels to the new "<year> Annual" format.
$ws.Range("B2").Value = "2024 Annual"
$ws.Range("C2").Value = "2023 Annual"
$ws.Range("D2").Value = "2022 Annual"
